$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Split "My mentor was Dr. Katie Morzinksi..." run into three runs
#    (no text content change, but matches canonical run layout).
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("My mentor was Dr. Katie Morzinksi", $true, $false, $false, $false, $false, $true, 1, $false, "My mentor was^&Dr. Katie Morzinksi", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Rewrite the "adaptive optics... image taking" description.
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("close-to-perfect astronomical image taking. The atmosphere actually distorts light from distant stars and planets, as the movement of hot and cool air above can cause disturbances to the light that reaches earth.", $true, $false, $false, $false, $false, $true, 1, $false, "improved surface astronomical image taking. When taking images of astronomical bodies from the surface of Earth, the atmosphere distorts light from distant stars and planets.", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("adaptive optics have emerged as a way to correct the problem of the atmosphere", $true, $false, $false, $false, $false, $true, 1, $false, "adaptive optics have emerged as one of the ways to correct the problem of the atmosphere", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Trim the Las Campagnas Observatory paragraph, removing the
#    "middle of nowhere" / light pollution discussion. This also
#    removes the _GoBack bookmark that previously lived inside that
#    deleted span (re-added later in its new location).
# ------------------------------------------------------------------
try {
    $bm = $d.Bookmarks("_GoBack")
    $bm.Delete()
} catch {
}

$r = $d.Content
$r.Find.Execute("in Chile at the Las Campagnas Observatory. As you can see here, this observatory is in the middle of nowhere, which is good for astronomical data collection. The higher up an observatory is, the less atmosphere there is to distort images. Also, there is not major sources of light pollution, such as cities, around this observatory, since it’s in the middle of the Atacama Desert. To the right", $true, $false, $false, $false, $false, $true, 1, $false, "in Chile at the Las Campagnas Observatory.^&To the right", 2) | Out-Null

# ------------------------------------------------------------------
# 4. Append new sentences after "see practical applications of this program."
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("see practical applications of this program.", $true, $false, $false, $false, $false, $true, 1, $false, "see practical applications of this program. I could experience a portion of the process of astronomical data collection, and view what happens to get the data that I was supposed to calibrate.", 2) | Out-Null

# ------------------------------------------------------------------
# 5. "being take." -> "being taken."
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("being take. It’s the detector equivalent", $true, $false, $false, $false, $false, $true, 1, $false, "being taken. It’s the detector equivalent", 2) | Out-Null

# ------------------------------------------------------------------
# 6. Expand "including numpy, matplotlib, and astropy." sentence.
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("including numpy, matplotlib, and astropy.", $true, $false, $false, $false, $false, $true, 1, $false, "including numpy to help with large arrays of values, matplotlib to graph, and astropy to work with the data files.", 2) | Out-Null

# ------------------------------------------------------------------
# 7. Insert a paragraph break (with a blank paragraph) between
#    "...compared to the data." and "In the end, it ended up being 4..."
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("deviation compared to the data. In the end, it ended up being 4", $true, $false, $false, $false, $false, $true, 1, $false, "deviation compared to the data. ^p^pIn the end, it ended up being 4", 2) | Out-Null

# ------------------------------------------------------------------
# 8. "a fairly straightforward process, but it actually took me" ->
#    "a straightforward process, but it took me"
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("it sounds like a fairly straightforward process, but it actually took me", $true, $false, $false, $false, $false, $true, 1, $false, "it sounds like a straightforward process, but it took me", 2) | Out-Null

# ------------------------------------------------------------------
# 9. Merge the "...to produce suitable corrections to the data." /
#    "Now, the important thing..." paragraphs into one (remove the
#    blank paragraph + paragraph break between them), with a new
#    space joining the two sentences.
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("months to produce suitable corrections to the data.^p^pNow, the important thing about these coefficients", $true, $false, $false, $false, $false, $true, 1, $false, "months to produce suitable corrections to the data. Now, the important thing about these coefficients", 2) | Out-Null

# ------------------------------------------------------------------
# 10. Merge "...For a better view of what I mean…" / "Let's talk
#     about these graphics." paragraphs into one (remove the blank
#     paragraph + paragraph break between them).
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("For a better view of what I mean…^p^pLet’s talk about these graphics.", $true, $false, $false, $false, $false, $true, 1, $false, "For a better view of what I mean…Let’s talk about these graphics.", 2) | Out-Null

# ------------------------------------------------------------------
# 11. Re-add the _GoBack bookmark inside "...represent the corrected
#     values..." (new location after the merge above).
# ------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("The black dotted lines represent the correct", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $target = $d.Range($r.End, $r.End)
    $d.Bookmarks.Add("_GoBack", $target) | Out-Null
}

# ------------------------------------------------------------------
# 12. Remove the trailing "Note the Airy disc pattern on the plot."
#     sentence.
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(" Note the Airy disc pattern on the plot.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ------------------------------------------------------------------
# 13. Drop the spell-check proofErr wrapper around "yous" (collapses
#     to a single run once the surrounding text is replaced as one).
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("The presentation has concluded, and thank yous go out to Katie Morzinski", $true, $false, $false, $false, $false, $true, 1, $false, "The presentation has concluded, and thank yous go out to Katie Morzinski", 2) | Out-Null

Write-Output "done"
